$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "drehen",
    "bitten",
    "führen",
    "schenken",
    "grüßen",
    "sorgen",
    "reizen",
    "treiben",
    "fahren",
    "lügen",
    "trauen",
    "werfen",
    "brauchen",
    "klingen",
    "kosten",
    "dringen",
    "schwingen",
    "töten",
    "geben",
    "hauen",
    "wirken",
    "stehlen",
    "jagen",
    "münzen",
    "zögern",
    "trennen",
    "schrecken",
    "gelten",
    "boxen",
    "malen",
    "schwören",
    "knarren",
    "greifen",
    "warnen",
    "spielen",
    "gründen",
    "bluten",
    "wehtun",
    "schwächen",
    "streichen",
    "schulden",
    "kichern",
    "flüchten",
    "formen",
    "werden",
    "spüren",
    "irren",
    "bellen",
    "ehren",
    "achten",
    "pfeifen",
    "wüten",
    "arten",
    "heben",
    "lockern",
    "schreiten",
    "sterben",
    "spinnen",
    "loben",
    "fallen",
    "liefern",
    "rasen",
    "wundern",
    "heilen",
    "stecken",
    "sprengen",
    "mauern",
    "sinken",
    "schmecken",
    "scheitern",
    "backen",
    "mögen",
    "kehren",
    "feiern",
    "äußern",
    "pflanzen",
    "enden",
    "zünden",
    "fällen",
    "flehen",
    "betteln",
    "quälen",
    "bergen",
    "räumen",
    "platzen",
    "bauen",
    "ärgern",
    "scheinen",
    "biegen",
    "graben",
    "freuen",
    "stammen",
    "machen",
    "planen",
    "tollen",
    "fangen",
    "tropfen",
    "lesen",
    "siegen",
    "kümmern",
    "zeigen",
    "ändern",
    "decken",
    "schlucken",
    "folgen",
    "hören",
    "zielen",
    "suchen",
    "runden",
    "sperren",
    "dienen",
    "saufen",
    "seufzen",
    "helfen",
    "jubeln",
    "wachsen",
    "rufen",
    "klettern",
    "fließen",
    "filmen",
    "heulen",
    "weichen",
    "erben",
    "altern",
    "sichern",
    "wenden",
    "liegen",
    "fischen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $words[$i]
}

Write-Output "Updated $($words.Length) word cells"
